$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.479.12"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.562.40"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "599.16"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "140.39"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "3.562.05"
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("E11").Value = "  -6.21%  "
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").Value = "4.166.90"
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "27.19"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.563.78"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "65.386.68"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "10.22"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").Value = "5.88"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "14.28"
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("D22").Value = "397.61"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "0.573"
$ws.Range("E23").Value = "  +4.66%  "
$ws.Range("D24").Value = "3.706.26"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("D25").Value = "74.68"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +10.16%  "
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +7.75%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "2.29"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "3.579.01"
$ws.Range("D33").Value = "23.97"
$ws.Range("E33").Value = "  +4.63%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "0.148"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").Value = "169.15"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").Value = "0.833"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").Value = "26.62"
$ws.Range("E43").Value = "  +15.63%  "
$ws.Range("D44").Value = "42.96"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("E46").Value = "  +4.97%  "
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +8.26%  "
$ws.Range("D49").Value = "2.452.06"
$ws.Range("E49").Value = "  +11.08%  "
$ws.Range("D50").Value = "6.83"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("E51").Value = "  +2.06%  "
